# Refresh cryptos list values (GitHub Actions scheduled scrape update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.746.14"
$ws.Range("E2").Value = "  +0.31%  "

# Row 3
$ws.Range("D3").Value = "1.603.09"
$ws.Range("E3").Value = "  +0.42%  "

# Row 4
$ws.Range("E4").Value = "  +0.15%  "

# Row 5
$ws.Range("E5").Value = "  +0.24%  "

# Row 6
$ws.Range("E6").Value = "  +0.41%  "

# Row 7
$ws.Range("E7").Value = "  +0.18%  "

# Row 8
$ws.Range("E8").Value = "  +0.28%  "

# Row 9
$ws.Range("E9").Value = "  +0.47%  "

# Row 10
$ws.Range("D10").Value = "'19.66"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.11%  "

# Row 11
$ws.Range("D11").Value = "'0.0847"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.74%  "

# Row 12
$ws.Range("D12").Value = "1.828.32"
$ws.Range("E12").Value = "  +0.40%  "

# Row 13
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.629.01"
$ws.Range("E13").Value = "  +0.79%  "

# Row 14
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'4.07"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.17%  "

# Row 15
$ws.Range("E15").Value = "  +0.50%  "

# Row 16
$ws.Range("E16").Value = "  +0.05%  "

# Row 17
$ws.Range("D17").Value = "0.0₃0742"
$ws.Range("E17").Value = "  +0.60%  "

# Row 18
$ws.Range("D18").Value = "'209.68"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.39%  "

# Row 19
$ws.Range("E19").Value = "  +0.19%  "

# Row 20
$ws.Range("D20").Value = "'7.15"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.26%  "

# Row 21
$ws.Range("E21").Value = "  +0.36%  "

# Row 22
$ws.Range("E22").Value = "  -4.81%  "

# Row 23
$ws.Range("E23").Value = "  +0.84%  "

# Row 24
$ws.Range("D24").Value = "'143.75"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.08%  "

# Row 25
$ws.Range("E25").Value = "  +0.13%  "

# Row 26
$ws.Range("E26").Value = "  -0.32%  "

# Row 27
$ws.Range("E27").Value = "  -0.19%  "

# Row 28
$ws.Range("D28").Value = "'15.38"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.60%  "

# Row 29
$ws.Range("D29").Value = "'0.0508"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.92%  "

# Row 30
$ws.Range("E30").Value = "  +0.20%  "

# Row 31
$ws.Range("D31").Value = "'3.28"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.49%  "

# Row 32
$ws.Range("D32").Value = "'2.97"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.86%  "

# Row 33
$ws.Range("D33").Value = "1.289.04"
$ws.Range("E33").Value = "  +0.09%  "

# Row 34
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").Value = "'2.48"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.23%  "

# Row 35
$ws.Range("B35").Value = "WEMIXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D35").Value = "'1.23"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +19.42%  "

# Row 36
$ws.Range("E36").Value = "  +0.44%  "

# Row 37
$ws.Range("D37").Value = "'0.590"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.76%  "

# Row 38
$ws.Range("E38").Value = "  -0.27%  "

# Row 39
$ws.Range("E39").Value = "  -0.18%  "

# Row 40
$ws.Range("D40").Value = "'5.45"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.00%  "

# Row 41
$ws.Range("E41").Value = "  -0.16%  "

# Row 42
$ws.Range("D42").Value = "'0.778"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.51%  "

# Row 43
$ws.Range("D43").Value = "'62.88"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.45%  "

# Row 44
$ws.Range("D44").Value = "1.739.89"
$ws.Range("E44").Value = "  +0.51%  "

# Row 45
$ws.Range("D45").Value = "'90.51"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.65%  "

# Row 46
$ws.Range("E46").Value = "  -0.04%  "

# Row 47
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "0.0₆0103"
$ws.Range("E47").Value = "  -3.25%  "

# Row 48
$ws.Range("B48").Value = "Algorand"
$ws.Range("C48").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D48").Value = "'0.102"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.12%  "

# Row 49
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "'0.0513"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.87%  "

# Row 50
$ws.Range("B50").Value = "Aptos"
$ws.Range("C50").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D50").Value = "'5.89"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +15.15%  "

# Row 51
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "'7.58"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.97%  "

Write-Host "Updated cryptos list"